$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 82 (existing rows 82-111 shift down to 84-113)
$ws.Rows(82).EntireRow.Insert()
$ws.Rows(82).EntireRow.Insert()

# New row 82: Cilantro, Primera
$ws.Range("A82").Value = 7
$ws.Range("B82").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C82").Value = "Ñuble"
$ws.Range("D82").Value = 44873
$ws.Range("E82").Value = 16
$ws.Range("F82").Value = 100112040
$ws.Range("G82").Value = "Cilantro"
$ws.Range("H82").Value = "Sin especificar"
$ws.Range("I82").Value = "Primera"
$ws.Range("J82").Value = 600
$ws.Range("K82").Value = 600
$ws.Range("L82").Value = 700
$ws.Range("M82").Value = 650
$ws.Range("N82").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O82").Value = "Provincia de Diguillín"
$ws.Range("P82").Value = 650
$ws.Range("Q82").Value = 1
$ws.Range("R82").Value = "Hortaliza"

# New row 83: Cilantro, Segunda
$ws.Range("A83").Value = 7
$ws.Range("B83").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C83").Value = "Ñuble"
$ws.Range("D83").Value = 44873
$ws.Range("E83").Value = 16
$ws.Range("F83").Value = 100112040
$ws.Range("G83").Value = "Cilantro"
$ws.Range("H83").Value = "Sin especificar"
$ws.Range("I83").Value = "Segunda"
$ws.Range("J83").Value = 500
$ws.Range("K83").Value = 500
$ws.Range("L83").Value = 500
$ws.Range("M83").Value = 500
$ws.Range("N83").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O83").Value = "Provincia de Diguillín"
$ws.Range("P83").Value = 500
$ws.Range("Q83").Value = 1
$ws.Range("R83").Value = "Hortaliza"
